$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the convective boundary condition temperatures in column D (rows 4-23):
# they were recorded in Fahrenheit instead of Kelvin, so shift every value
# up by 300 to correct them.
for ($r = 4; $r -le 23; $r++) {
    $cell = $ws.Range("D$r")
    $cell.Value2 = $cell.Value2 + 300
}

# Add a new (currently empty) column L for each data row, formatted with the
# "0.0" number format and centered alignment, ready for Transient variable
# names to be filled in later.
for ($r = 4; $r -le 23; $r++) {
    $lcell = $ws.Range("L$r")
    $lcell.NumberFormat = "0.0"
    $lcell.HorizontalAlignment = -4108
    $lcell.VerticalAlignment = -4108
}

# Restore the active selection to where the author left off editing.
$ws.Range("L17").Select()
